$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '60.594.61'
$ws.Range("E2").Value = '  +5.99%  '
$ws.Range("D3").Value = '2.634.36'
$ws.Range("E3").Value = '  +8.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.96'
$ws.Range("E5").Value = '  +3.48%  '
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  -5.11%  '
$ws.Range("D9").Value = '2.633.96'
$ws.Range("E9").Value = '  +7.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.46'
$ws.Range("E10").Value = '  +5.01%  '
$ws.Range("E11").Value = '  +4.01%  '
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("E13").Value = '  +0.89%  '
$ws.Range("D14").Value = '3.064.85'
$ws.Range("E14").Value = '  +7.00%  '
$ws.Range("D15").Value = '60.661.75'
$ws.Range("E15").Value = '  +6.02%  '
$ws.Range("E16").Value = '  +4.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  +4.58%  '
$ws.Range("D18").Value = '2.620.04'
$ws.Range("E18").Value = '  +7.16%  '
$ws.Range("E19").Value = '  +2.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.27'
$ws.Range("E20").Value = '  +5.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.41'
$ws.Range("E21").Value = '  +3.97%  '
$ws.Range("E22").Value = '  +3.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.32'
$ws.Range("E24").Value = '  +4.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.422'
$ws.Range("E25").Value = '  +5.12%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.165'
$ws.Range("E26").Value = '  +2.79%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0853'
$ws.Range("E28").Value = '  +8.78%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.54'
$ws.Range("E29").Value = '  +3.47%  '
$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.72'
$ws.Range("E31").Value = '  +3.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.43'
$ws.Range("E32").Value = '  +3.91%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.57'
$ws.Range("E33").Value = '  +3.08%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.72'
$ws.Range("E34").Value = '  +7.79%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.01'
$ws.Range("E35").Value = '  +5.87%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.20'
$ws.Range("E36").Value = '  +4.53%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '307.18'
$ws.Range("E37").Value = '  +7.05%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.846'
$ws.Range("E38").Value = '  +3.59%  '
$ws.Range("E39").Value = '  +7.31%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.77'
$ws.Range("E40").Value = '  +6.95%  '
$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.827'
$ws.Range("E41").Value = '  +26.01%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.56'
$ws.Range("E42").Value = '  +4.61%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.628'
$ws.Range("E43").Value = '  +4.69%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0570'
$ws.Range("E44").Value = '  +7.10%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.100'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.82'
$ws.Range("E47").Value = '  +12.57%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.88'
$ws.Range("E48").Value = '  +8.36%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0235'
$ws.Range("E49").Value = '  +3.42%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.046.69'
$ws.Range("E50").Value = '  +7.57%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.27'
$ws.Range("E51").Value = '  +0.46%  '
